$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 239
$ws1.Range("F3").Value = 264
$ws1.Range("F4").Value = 277
$ws1.Range("F5").Value = 820
$ws1.Range("F7").Value = 6574
$ws1.Range("F8").Value = 52
$ws1.Range("F15").Value = 210
$ws1.Range("F16").Value = 535
$ws1.Range("F17").Value = 54

# Sheet "全部类型" (All Types) - same data duplicated
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 239
$ws4.Range("F3").Value = 264
$ws4.Range("F4").Value = 277
$ws4.Range("F5").Value = 820
$ws4.Range("F7").Value = 6574
$ws4.Range("F8").Value = 52
$ws4.Range("F15").Value = 210
$ws4.Range("F16").Value = 535
$ws4.Range("F17").Value = 54
